# #5: cash & deposit done
# Rebuild the "存款" (deposit) sheet: add bank / deposit_type / currency
# columns up front (B/C/D), fix the header row to use real field names,
# and append the normalized metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- 1. Extend formatting from the existing columns onto the new ones ---
# Header row (style 1: bold + border + centered) -- tile B1:C1 across G1:M1
$ws.Range("B1:C1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

# Data rows (style 2 / default look) -- tile B2:C2 across G2:M20
$ws.Range("B2:C2").Copy()
$ws.Range("G2:M20").PasteSpecial(-4122)

# --- 2. Header row values ---
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- 3. Data rows (row 2 = index 58 ... row 20 = index 77) ---
$rows = @(
  @{ R=2;  A=58; B="合作金庫商業銀行古亭分行";         C="綜合存款";     D="新臺幣"; E="陳學聖"; F=716714;  M=58 },
  @{ R=3;  A=59; B="合作金庫商業銀行永吉分行";         C="活期儲蓄存款"; D="新臺幣"; E="陳學聖"; F=68309;   M=59 },
  @{ R=4;  A=60; B="台北富邦商業銀行古亭分行";         C="活期儲蓄存款"; D="新臺幣"; E="陳學聖"; F=24859;   M=60 },
  @{ R=5;  A=61; B="大眾商業銀行前鎮分行";             C="活期儲蓄存款"; D="新臺幣"; E="陳學聖"; F=203;     M=61 },
  @{ R=6;  A=62; B="臺灣土地銀行古亭分行";             C="活期儲蓄存款"; D="新臺幣"; E="陳學聖"; F=65;      M=62 },
  @{ R=7;  A=63; B="中華郵政股份有限公司";             C="活期存款";     D="新臺幣"; E="陳學聖"; F=145595;  M=63 },
  @{ R=8;  A=65; B="中華郵政股份有限公司新店青潭郵局"; C="活期存款";     D="新臺幣"; E="陳學聖"; F=6415914; M=65 },
  @{ R=9;  A=66; B="臺灣中小企業銀行世貿分行";         C="活期儲蓄存款"; D="新臺幣"; E="陳學聖"; F=143553;  M=66 },
  @{ R=10; A=67; B="中華郵政股份有限公司新店青潭郵局"; C="活期存款";     D="新臺幣"; E="梁寒衣"; F=1994108; M=67 },
  @{ R=11; A=68; B="第一商業銀行吉林簡易型分行";       C="活期儲蓄存款"; D="新臺幣"; E="梁寒衣"; F=215081;  M=68 },
  @{ R=12; A=69; B="台北富邦商業銀行";                 C="支票存款";     D="新臺幣"; E="陳學聖"; F=327;     M=69 },
  @{ R=13; A=70; B="中華郵政股份有限公司";             C="活期存款";     D="新臺幣"; E="陳學聖"; F=876;     M=70 },
  @{ R=14; A=71; B="聯邦商業銀行";                     C="活期存款";     D="新臺幣"; E="陳學聖"; F=3294;    M=71 },
  @{ R=15; A=72; B="聯邦商業銀行";                     C="活期儲蓄存款"; D="新臺幣"; E="陳學聖"; F=54;      M=72 },
  @{ R=16; A=73; B="臺灣銀行";                         C="活期儲蓄存款"; D="新臺幣"; E="陳學聖"; F=1051;    M=73 },
  @{ R=17; A=74; B="彰化商業銀行古亭分行";             C="活期存款";     D="新臺幣"; E="陳學聖"; F=1104;    M=74 },
  @{ R=18; A=75; B="華南商業銀行";                     C="活期儲蓄存款"; D="新臺幣"; E="陳學聖"; F=22697;   M=75 },
  @{ R=19; A=76; B="合作金庫商業銀行新店分行";         C="活期儲蓄存款"; D="新臺幣"; E="梁寒衣"; F=1694;    M=76 },
  @{ R=20; A=77; B="合作金庫商業銀行古亭分行";         C="活期存款";     D="美金";   E="陳學聖"; F=2183.66; M=77 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value  = $row.A   # A: index (legacy leading column)
    $ws.Cells.Item($r, 2).Value  = $row.B   # B: bank
    $ws.Cells.Item($r, 3).Value  = $row.C   # C: deposit_type
    $ws.Cells.Item($r, 4).Value  = $row.D   # D: currency
    $ws.Cells.Item($r, 5).Value  = $row.E   # E: owner
    $ws.Cells.Item($r, 6).Value  = $row.F   # F: total
    $ws.Cells.Item($r, 7).Value  = "deposit"      # G: property_category
    $ws.Cells.Item($r, 8).Value  = "normal"       # H: category
    $ws.Cells.Item($r, 9).Value  = "2012-04-25"   # I: date
    $ws.Cells.Item($r, 10).Value = "陳學聖"        # J: legislator_name
    $ws.Cells.Item($r, 11).Value = 840            # K: legislator_id
    $ws.Cells.Item($r, 12).Value = "tmpfd9c1"     # L: source_file
    $ws.Cells.Item($r, 13).Value = $row.M         # M: index
}
